$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with the new match data (league, date, time, teams, odds).
# Force Date/Time columns (B, C) to Text format first so values like
# '2025-12-23' are kept as literal strings instead of being auto-converted
# to date serials, matching the source data (inline strings).
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = 'Friendly Matches'
$ws.Cells.Item(2,2).Value = '2025-12-23'
$ws.Cells.Item(2,3).Value = '16:00:00'
$ws.Cells.Item(2,4).Value = 'Serra Branca EC'
$ws.Cells.Item(2,5).Value = 'Maguary'
$ws.Cells.Item(2,6).Value = 1.04
$ws.Cells.Item(2,7).Value = 980
$ws.Cells.Item(2,8).Value = 1.09
$ws.Cells.Item(2,9).Value = 1000
$ws.Cells.Item(2,10).Value = 3.5
$ws.Cells.Item(2,11).Value = 3.8
$ws.Cells.Item(2,12).Value = 1.03
$ws.Cells.Item(2,13).Value = 1.07
$ws.Cells.Item(2,14).Value = 3.95
$ws.Cells.Item(2,15).Value = 1.04
$ws.Cells.Item(2,16).Value = 1.04
$ws.Cells.Item(2,17).Value = 1.05
$ws.Cells.Item(2,18).Value = 1.05
$ws.Cells.Item(2,19).Value = 1.02
$ws.Cells.Item(2,20).Value = 1.63
$ws.Cells.Item(2,21).Value = 2.08
$ws.Cells.Item(2,22).Value = 1.02
$ws.Cells.Item(2,23).Value = 1.02
$ws.Cells.Item(2,24).Value = 990
$ws.Cells.Item(2,25).Value = 990
$ws.Cells.Item(2,26).Value = 980
$ws.Cells.Item(2,27).Value = 120
$ws.Cells.Item(2,28).Value = 990
$ws.Cells.Item(2,29).Value = 990
$ws.Cells.Item(2,30).Value = 990
$ws.Cells.Item(2,31).Value = 980
$ws.Cells.Item(2,32).Value = 980
$ws.Cells.Item(2,33).Value = 990
$ws.Cells.Item(2,34).Value = 990
$ws.Cells.Item(2,35).Value = 980
$ws.Cells.Item(2,36).Value = 980
$ws.Cells.Item(2,37).Value = 980
$ws.Cells.Item(2,38).Value = 980
$ws.Cells.Item(2,39).Value = 140
$ws.Cells.Item(2,40).Value = 980
$ws.Cells.Item(2,41).Value = 980

$ws.Cells.Item(3,2).NumberFormat = "@"
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = 'Portuguese Primeira Liga'
$ws.Cells.Item(3,2).Value = '2025-12-23'
$ws.Cells.Item(3,3).Value = '17:45:00'
$ws.Cells.Item(3,4).Value = 'Guimaraes'
$ws.Cells.Item(3,5).Value = 'Sporting Lisbon'
$ws.Cells.Item(3,6).Value = 27
$ws.Cells.Item(3,7).Value = 29
$ws.Cells.Item(3,8).Value = 1.22
$ws.Cells.Item(3,9).Value = 1.23
$ws.Cells.Item(3,10).Value = 6.4
$ws.Cells.Item(3,11).Value = 6.8
$ws.Cells.Item(3,12).Value = 0
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,14).Value = 0
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = 0
$ws.Cells.Item(3,17).Value = 0
$ws.Cells.Item(3,18).Value = 3.85
$ws.Cells.Item(3,19).Value = 1.33
$ws.Cells.Item(3,20).Value = 0
$ws.Cells.Item(3,21).Value = 0
$ws.Cells.Item(3,22).Value = 5.3
$ws.Cells.Item(3,23).Value = 1.03
$ws.Cells.Item(3,24).Value = 1000
$ws.Cells.Item(3,25).Value = 1000
$ws.Cells.Item(3,26).Value = 1000
$ws.Cells.Item(3,27).Value = 1000
$ws.Cells.Item(3,28).Value = 1000
$ws.Cells.Item(3,29).Value = 1000
$ws.Cells.Item(3,30).Value = 4.2
$ws.Cells.Item(3,31).Value = 4.5
$ws.Cells.Item(3,32).Value = 1000
$ws.Cells.Item(3,33).Value = 1000
$ws.Cells.Item(3,34).Value = 8.6
$ws.Cells.Item(3,35).Value = 9.2
$ws.Cells.Item(3,36).Value = 1000
$ws.Cells.Item(3,37).Value = 1000
$ws.Cells.Item(3,38).Value = 46
$ws.Cells.Item(3,39).Value = 38
$ws.Cells.Item(3,40).Value = 160
$ws.Cells.Item(3,41).Value = 4.6

$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = 'Friendly Matches'
$ws.Cells.Item(4,2).Value = '2025-12-23'
$ws.Cells.Item(4,3).Value = '18:00:00'
$ws.Cells.Item(4,4).Value = 'Necaxa'
$ws.Cells.Item(4,5).Value = 'Atletico San Luis'
$ws.Cells.Item(4,6).Value = 1.09
$ws.Cells.Item(4,7).Value = 1000
$ws.Cells.Item(4,8).Value = 1.09
$ws.Cells.Item(4,9).Value = 1000
$ws.Cells.Item(4,10).Value = 1.03
$ws.Cells.Item(4,11).Value = 1000
$ws.Cells.Item(4,12).Value = 1.03
$ws.Cells.Item(4,13).Value = 1.07
$ws.Cells.Item(4,14).Value = 3.7
$ws.Cells.Item(4,15).Value = 1.02
$ws.Cells.Item(4,16).Value = 1.25
$ws.Cells.Item(4,17).Value = 1.03
$ws.Cells.Item(4,18).Value = 1.18
$ws.Cells.Item(4,19).Value = 1.01
$ws.Cells.Item(4,20).Value = 1.65
$ws.Cells.Item(4,21).Value = 2.06
$ws.Cells.Item(4,22).Value = 1.03
$ws.Cells.Item(4,23).Value = 1.03
$ws.Cells.Item(4,24).Value = 1000
$ws.Cells.Item(4,25).Value = 990
$ws.Cells.Item(4,26).Value = 32
$ws.Cells.Item(4,27).Value = 1000
$ws.Cells.Item(4,28).Value = 990
$ws.Cells.Item(4,29).Value = 990
$ws.Cells.Item(4,30).Value = 990
$ws.Cells.Item(4,31).Value = 50
$ws.Cells.Item(4,32).Value = 1000
$ws.Cells.Item(4,33).Value = 990
$ws.Cells.Item(4,34).Value = 1000
$ws.Cells.Item(4,35).Value = 1000
$ws.Cells.Item(4,36).Value = 34
$ws.Cells.Item(4,37).Value = 32
$ws.Cells.Item(4,38).Value = 80
$ws.Cells.Item(4,39).Value = 1000
$ws.Cells.Item(4,40).Value = 970
$ws.Cells.Item(4,41).Value = 50

$ws.Cells.Item(5,2).NumberFormat = "@"
$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = 'Honduras Liga Nacional'
$ws.Cells.Item(5,2).Value = '2025-12-23'
$ws.Cells.Item(5,3).Value = '22:00:00'
$ws.Cells.Item(5,4).Value = 'Real Espana'
$ws.Cells.Item(5,5).Value = 'CD Motagua'
$ws.Cells.Item(5,6).Value = 1.71
$ws.Cells.Item(5,7).Value = 1.8
$ws.Cells.Item(5,8).Value = 5.1
$ws.Cells.Item(5,9).Value = 6.2
$ws.Cells.Item(5,10).Value = 3.95
$ws.Cells.Item(5,11).Value = 4.3
$ws.Cells.Item(5,12).Value = 1.41
$ws.Cells.Item(5,13).Value = 1.07
$ws.Cells.Item(5,14).Value = 3.55
$ws.Cells.Item(5,15).Value = 1.32
$ws.Cells.Item(5,16).Value = 1.88
$ws.Cells.Item(5,17).Value = 1.95
$ws.Cells.Item(5,18).Value = 1.34
$ws.Cells.Item(5,19).Value = 3.5
$ws.Cells.Item(5,20).Value = 1.85
$ws.Cells.Item(5,21).Value = 1.92
$ws.Cells.Item(5,22).Value = 1.2
$ws.Cells.Item(5,23).Value = 2.26
$ws.Cells.Item(5,24).Value = 14.5
$ws.Cells.Item(5,25).Value = 19
$ws.Cells.Item(5,26).Value = 48
$ws.Cells.Item(5,27).Value = 160
$ws.Cells.Item(5,28).Value = 8.2
$ws.Cells.Item(5,29).Value = 10
$ws.Cells.Item(5,30).Value = 23
$ws.Cells.Item(5,31).Value = 85
$ws.Cells.Item(5,32).Value = 10.5
$ws.Cells.Item(5,33).Value = 10.5
$ws.Cells.Item(5,34).Value = 23
$ws.Cells.Item(5,35).Value = 85
$ws.Cells.Item(5,36).Value = 18.5
$ws.Cells.Item(5,37).Value = 19
$ws.Cells.Item(5,38).Value = 38
$ws.Cells.Item(5,39).Value = 140
$ws.Cells.Item(5,40).Value = 12
$ws.Cells.Item(5,41).Value = 100

# Remove the now-obsolete rows 6 and 7 (their matches were consolidated
# into rows 2-5 above), shrinking the used range to A1:AO5.
$ws.Rows("6:7").Delete()
